# 5000_attendance.xlsx
# Employee switched from "Temporary" (code 70001004) to "Muhammad Usman" (code 5000).
# The attendance rows are replaced with six "Applied For Leave" rows (Oct 02-07 2023).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column width changes --------------------------------------------------
# Former merged "B:C" width-18 run is split: B stays 18, C widens to fit names.
$ws.Columns("C").ColumnWidth = 18.1666666666667
# Status column widens to fit "Applied For Leave".
$ws.Columns("J").ColumnWidth = 21.1666666666667

# ---- Build rows 5-7 (new rows) by cloning row 2's cell formatting ---------
$ws.Range("A2:L2").Copy()
$ws.Range("A5:L5").PasteSpecial(-4122)
$ws.Range("A2:L2").Copy()
$ws.Range("A6:L6").PasteSpecial(-4122)
$ws.Range("A2:L2").Copy()
$ws.Range("A7:L7").PasteSpecial(-4122)

# Drop the Start/End Time + Break columns entirely for every data row (E:H) -
# "Applied For Leave" entries carry no punch times.
$ws.Range("E2:H2").Clear()
$ws.Range("E3:H3").Clear()
$ws.Range("E4:H4").Clear()
$ws.Range("E5:H5").Clear()
$ws.Range("E6:H6").Clear()
$ws.Range("E7:H7").Clear()

# Status cells no longer get the green "Paid" highlight - reformat them with
# the plain body style (same style already used by column L).
$ws.Range("L2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("J5").PasteSpecial(-4122)
$ws.Range("J6").PasteSpecial(-4122)
$ws.Range("J7").PasteSpecial(-4122)

# ---- Cell values ------------------------------------------------------------
$dates = @("Sat Oct 07 2023", "Fri Oct 06 2023", "Thu Oct 05 2023", "Wed Oct 04 2023", "Tue Oct 03 2023", "Mon Oct 02 2023")

for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    $ws.Range("A$r").Value = $i + 1
    $ws.Range("B$r").Value = 5000
    $ws.Range("C$r").Value = "Muhammad Usman"
    $ws.Range("D$r").Value = $dates[$i]
    $ws.Range("I$r").Value = "---"
    $ws.Range("J$r").Value = "Applied For Leave"
    $ws.Range("K$r").Value = "'"
    $ws.Range("L$r").Value = "'"
}
